$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.645.54"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "2.500.14"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("D5").Value = "'590.38"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").Value = "'173.33"
$ws.Range("E6").Value = "  -1.95%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.514"
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").Value = "2.500.51"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  +3.93%  "
$ws.Range("D12").Value = "'4.97"
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").Value = "'0.333"
$ws.Range("E13").Value = "  -2.62%  "
$ws.Range("D14").Value = "2.951.84"
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("D15").Value = "'25.49"
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("D16").Value = "68.567.07"
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("D17").Value = "'0.0000171"
$ws.Range("E17").Value = "  -1.02%  "
$ws.Range("D18").Value = "2.511.86"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("D19").Value = "'357.78"
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("D20").Value = "'7.46"
$ws.Range("E20").Value = "  -1.16%  "
$ws.Range("D21").Value = "'10.81"
$ws.Range("E21").Value = "  -2.87%  "
$ws.Range("D22").Value = "'3.99"
$ws.Range("E22").Value = "  -2.66%  "
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").Value = "'70.05"
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("D25").Value = "'4.11"
$ws.Range("E25").Value = "  -4.84%  "
$ws.Range("D26").Value = "'8.80"
$ws.Range("E26").Value = "  -4.30%  "
$ws.Range("D27").Value = "2.631.14"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("D28").Value = "'1.62"
$ws.Range("E28").Value = "  -8.12%  "
$ws.Range("D29").Value = "'0.997"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "'503.70"
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("D31").Value = "0.0₃0866"
$ws.Range("E31").Value = "  -5.25%  "
$ws.Range("D32").Value = "'7.64"
$ws.Range("E32").Value = "  -2.71%  "
$ws.Range("D33").Value = "'1.76"
$ws.Range("E33").Value = "  -1.36%  "
$ws.Range("D34").Value = "'1.19"
$ws.Range("E34").Value = "  -5.73%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").Value = "'161.92"
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("D37").Value = "'0.117"
$ws.Range("E37").Value = "  -5.44%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").Value = "'18.63"
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").Value = "'1.29"
$ws.Range("E41").Value = "  -4.03%  "
$ws.Range("D42").Value = "'1.68"
$ws.Range("E42").Value = "  -3.89%  "
$ws.Range("D43").Value = "'4.67"
$ws.Range("E43").Value = "  -4.33%  "
$ws.Range("D44").Value = "'0.314"
$ws.Range("E44").Value = "  -5.09%  "
$ws.Range("D45").Value = "'2.27"
$ws.Range("E45").Value = "  -6.28%  "
$ws.Range("D46").Value = "'148.87"
$ws.Range("E46").Value = "  +1.78%  "
$ws.Range("D47").Value = "'3.51"
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("D48").Value = "'0.506"
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("D49").Value = "'0.0732"
$ws.Range("E49").Value = "  -2.02%  "
$ws.Range("D50").Value = "'1.55"
$ws.Range("E50").Value = "  -3.18%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.573"
$ws.Range("E51").Value = "  -2.46%  "
